$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: insert "accuracy" after "model", shift others right,
# and add "smote" as the new last header in F1.
$ws.Range("A1").Value = "model"
$ws.Range("B1").Value = "accuracy"
$ws.Range("C1").Value = "sensitivity"
$ws.Range("D1").Value = "specificity"
$ws.Range("E1").Value = "precision"
$ws.Range("F1").Value = "smote"

# F1 needs the same header styling (bold, centered, bordered) as the rest
# of the header row.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "smote"

# Row 2: gradient_boosting results
$ws.Range("A2").Value = "gradient_boosting"
$ws.Range("B2").Value = 0.9970788704965921
$ws.Range("C2").Value = 0.9655172413793104
$ws.Range("D2").Value = 0.9979959919839679
$ws.Range("E2").Value = 0.9333333333333333
$ws.Range("F2").Value = $false

# Row 3: logistic_regression results
$ws.Range("A3").Value = "logistic_regression"
$ws.Range("B3").Value = 0.9863680623174295
$ws.Range("C3").Value = 0.5172413793103449
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = $false

# Remove the old fourth row (table now only has 2 data rows)
$ws.Range("A4:F4").Delete()
